$d = $word.ActiveDocument

$replacements = @(
    @{old="729×2=1458"; new="513×4=2052"},
    @{old="269×7=1883"; new="899×7=6293"},
    @{old="569×9=5121"; new="394×9=3546"},
    @{old="504×6=3024"; new="790×2=1580"},
    @{old="646×4=2584"; new="707×8=5656"},
    @{old="238×6=1428"; new="802×5=4010"},
    @{old="581×8=4648"; new="528×6=3168"},
    @{old="871×2=1742"; new="438×2=876"},
    @{old="821×2=1642"; new="978×5=4890"},
    @{old="944×5=4720"; new="123×4=492"},
    @{old="643×8=5144"; new="702×9=6318"},
    @{old="815×7=5705"; new="299×4=1196"},
    @{old="987×9=8883"; new="200×3=600"},
    @{old="438×5=2190"; new="337×3=1011"},
    @{old="262×4=1048"; new="651×7=4557"},
    @{old="441×4=1764"; new="715×8=5720"},
    @{old="394×4=1576"; new="510×5=2550"},
    @{old="997×9=8973"; new="517×3=1551"},
    @{old="790×6=4740"; new="518×6=3108"},
    @{old="757×7=5299"; new="701×7=4907"},
    @{old="681×8=5448"; new="539×3=1617"},
    @{old="690×3=2070"; new="308×7=2156"},
    @{old="321×5=1605"; new="916×3=2748"},
    @{old="822×6=4932"; new="519×5=2595"},
    @{old="122×2=244"; new="816×4=3264"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
